$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 9 (only 8 data rows remain: rows 2-8, header row 1)
$ws.Rows(9).Delete()

# Clear column E entirely (avg_popularity values removed; header stays same text but column E is now unused -> but per diff, E1 header text stays "avg_popularity")
# Actually E1 header remains avg_popularity; only E2:E8 data removed
$ws.Range("E2").ClearContents()
$ws.Range("E3").ClearContents()
$ws.Range("E4").ClearContents()
$ws.Range("E5").ClearContents()
$ws.Range("E6").ClearContents()
$ws.Range("E7").ClearContents()
$ws.Range("E8").ClearContents()

# Update data cells per target
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 2023
$ws.Range("D2").Value = 230
$ws.Range("F2").Value = 114.722479185939
$ws.Range("G2").Value = 8
$ws.Range("H2").Value = "Orlando"
$ws.Range("I2").Value = "Boston"
$ws.Range("J2").Value = 0.5171717171717172
$ws.Range("K2").Value = 98.42509250693803
$ws.Range("L2").Value = 115.8052728954672
$ws.Range("M2").Value = 114.4395004625347
$ws.Range("N2").Value = 77.55603607770584
$ws.Range("O2").Value = 0.4198704902867715
$ws.Range("P2").Value = 0.5908640148011101
$ws.Range("Q2").Value = 0.2783799722479186
$ws.Range("R2").Value = 12.37654949121184
$ws.Range("S2").Value = 11.5350138760407
$ws.Range("T2").Value = 0.2105568917668825
$ws.Range("U2").Value = 1.005455558158974
$ws.Range("V2").Value = 1.038430802510143
$ws.Range("W2").Value = 10.46283211763107
$ws.Range("X2").Value = 0.557123034227567
$ws.Range("Y2").Value = 40.5
$ws.Range("Z2").Value = 76.6
$ws.Range("AA2").Value = 0.4936169592900065

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 2023
$ws.Range("D3").Value = 236
$ws.Range("F3").Value = 112.2164855072464
$ws.Range("G3").Value = 11
$ws.Range("H3").Value = "Detroit"
$ws.Range("I3").Value = "Milwaukee"
$ws.Range("J3").Value = 0.5372340425531915
$ws.Range("K3").Value = 98.9381793478261
$ws.Range("L3").Value = 112.646240942029
$ws.Range("M3").Value = 115.878125
$ws.Range("N3").Value = 76.05253623188406
$ws.Range("O3").Value = 0.4087871376811594
$ws.Range("P3").Value = 0.5646467391304347
$ws.Range("Q3").Value = 0.2813088768115942
$ws.Range("R3").Value = 12.52803442028986
$ws.Range("S3").Value = 11.2513134057971
$ws.Range("T3").Value = 0.2089909420289855
$ws.Range("U3").Value = 0.9834924233763924
$ws.Range("V3").Value = 1.01876925814006
$ws.Range("W3").Value = 11.46413181898958
$ws.Range("X3").Value = 0.4402173913043478
$ws.Range("Y3").Value = 41.5
$ws.Range("Z3").Value = 75.15
$ws.Range("AA3").Value = 0.5349346164410206

$ws.Range("A4").Value = 2
$ws.Range("B4").Value = 2023
$ws.Range("D4").Value = 239.5
$ws.Range("F4").Value = 115.3787234042553
$ws.Range("G4").Value = 1
$ws.Range("H4").Value = "Chicago"
$ws.Range("I4").Value = "Atlanta"
$ws.Range("J4").Value = 0.531678486997636
$ws.Range("K4").Value = 99.58307328605201
$ws.Range("L4").Value = 114.6597163120567
$ws.Range("M4").Value = 114.905011820331
$ws.Range("N4").Value = 76.71513002364065
$ws.Range("O4").Value = 0.3341248226950354
$ws.Range("P4").Value = 0.5818765957446808
$ws.Range("Q4").Value = 0.248455791962175
$ws.Range("R4").Value = 11.37695035460993
$ws.Range("S4").Value = 12.54451536643026
$ws.Range("T4").Value = 0.2077384160756502
$ws.Range("U4").Value = 1.011207041229232
$ws.Range("V4").Value = 1.09050101759209
$ws.Range("W4").Value = 10.71066028731968
$ws.Range("X4").Value = 0.4886524822695035
$ws.Range("Y4").Value = 44
$ws.Range("Z4").Value = 76.6
$ws.Range("AA4").Value = 0.478910851199879

$ws.Range("A5").Value = 3
$ws.Range("B5").Value = 2023
$ws.Range("D5").Value = 234
$ws.Range("F5").Value = 112.4932065217391
$ws.Range("G5").Value = 5.5
$ws.Range("H5").Value = "Houston"
$ws.Range("I5").Value = "Minnesota"
$ws.Range("J5").Value = 0.4895833333333334
$ws.Range("K5").Value = 99.61326992753621
$ws.Range("L5").Value = 112.5564764492754
$ws.Range("M5").Value = 116.633786231884
$ws.Range("N5").Value = 74.7070652173913
$ws.Range("O5").Value = 0.3866417572463768
$ws.Range("P5").Value = 0.5744692028985507
$ws.Range("Q5").Value = 0.2839759963768115
$ws.Range("R5").Value = 13.71272644927536
$ws.Range("S5").Value = 12.50629528985507
$ws.Range("T5").Value = 0.2192010869565217
$ws.Range("U5").Value = 0.9859176732843044
$ws.Range("V5").Value = 1.055318573668062
$ws.Range("W5").Value = 10.84714287100455
$ws.Range("X5").Value = 0.358695652173913
$ws.Range("Y5").Value = 36.5
$ws.Range("Z5").Value = 74.7
$ws.Range("AA5").Value = 0.5008239892537384

$ws.Range("A6").Value = 4
$ws.Range("B6").Value = 2023
$ws.Range("D6").Value = 239
$ws.Range("F6").Value = 115.0423360833695
$ws.Range("G6").Value = 8
$ws.Range("H6").Value = "Utah"
$ws.Range("I6").Value = "Charlotte"
$ws.Range("J6").Value = 0.5472074468085106
$ws.Range("K6").Value = 99.76076856274423
$ws.Range("L6").Value = 114.3002605297438
$ws.Range("M6").Value = 117.0547546678246
$ws.Range("N6").Value = 75.2677811550152
$ws.Range("O6").Value = 0.402405557967868
$ws.Range("P6").Value = 0.568684759009987
$ws.Range("Q6").Value = 0.2626880156317847
$ws.Range("R6").Value = 12.07757273122015
$ws.Range("S6").Value = 12.12581415544942
$ws.Range("T6").Value = 0.2098382544507164
$ws.Range("U6").Value = 1.008258861379224
$ws.Range("V6").Value = 1.043698859086146
$ws.Range("W6").Value = 10.76685529504964
$ws.Range("X6").Value = 0.383195831524099
$ws.Range("Y6").Value = 29
$ws.Range("Z6").Value = 74.94999999999999
$ws.Range("AA6").Value = 0.5169201007426152

$ws.Range("A7").Value = 5
$ws.Range("B7").Value = 2023
$ws.Range("D7").Value = 241
$ws.Range("F7").Value = 112.804347826087
$ws.Range("G7").Value = 8
$ws.Range("H7").Value = "Portland"
$ws.Range("I7").Value = "SanAntonio"
$ws.Range("J7").Value = 0.5227272727272727
$ws.Range("K7").Value = 98.79999999999998
$ws.Range("L7").Value = 113.9260869565217
$ws.Range("M7").Value = 118.2184782608696
$ws.Range("N7").Value = 75.9891304347826
$ws.Range("O7").Value = 0.3754347826086956
$ws.Range("P7").Value = 0.5772934782608695
$ws.Range("Q7").Value = 0.2719891304347826
$ws.Range("R7").Value = 13.15543478260869
$ws.Range("S7").Value = 11.91413043478261
$ws.Range("T7").Value = 0.2077445652173913
$ws.Range("U7").Value = 0.988644590938536
$ws.Range("V7").Value = 0.9914098480645845
$ws.Range("W7").Value = 11.52216595171919
$ws.Range("X7").Value = 0.3804347826086957
$ws.Range("Y7").Value = 31
$ws.Range("Z7").Value = 74.15
$ws.Range("AA7").Value = 0.4977575470846392

$ws.Range("A8").Value = 6
$ws.Range("B8").Value = 2023
$ws.Range("D8").Value = 246.5
$ws.Range("F8").Value = 118.7065217391304
$ws.Range("G8").Value = 1
$ws.Range("H8").Value = "Sacramento"
$ws.Range("I8").Value = "Memphis"
$ws.Range("J8").Value = 0.5
$ws.Range("K8").Value = 100.8882367149758
$ws.Range("L8").Value = 117.3172946859904
$ws.Range("M8").Value = 113.2206280193237
$ws.Range("N8").Value = 77.02420289855073
$ws.Range("O8").Value = 0.38692922705314
$ws.Range("P8").Value = 0.5855654589371981
$ws.Range("Q8").Value = 0.2835183574879227
$ws.Range("R8").Value = 11.81096618357488
$ws.Range("S8").Value = 12.19275362318841
$ws.Range("T8").Value = 0.2119192028985507
$ws.Range("U8").Value = 1.040372670807453
$ws.Range("V8").Value = 0.9925325042687567
$ws.Range("W8").Value = 11.72748061756648
$ws.Range("X8").Value = 0.6258454106280193
$ws.Range("Y8").Value = 42
$ws.Range("Z8").Value = 75.55000000000001
$ws.Range("AA8").Value = 0.4889823474744299

$ws.Range("AA1").Value = "calc_over_prob"
